$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("courses")

# Update the department name for the course row from the old faculty name
# to the new shorter department label "Community Services".
$ws.Range("C2").Value = "Community Services"
